$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 19
$ws1.Range("F5").Value = 7731
$ws1.Range("F7").Value = 114
$ws1.Range("F8").Value = 2105
$ws1.Range("F9").Value = 8513
$ws1.Range("F10").Value = 11
$ws1.Range("F12").Value = 84
$ws1.Range("F13").Value = 5711
$ws1.Range("F15").Value = 2661
$ws1.Range("F16").Value = 1169
$ws1.Range("F17").Value = 4600
$ws1.Range("F18").Value = 353
$ws1.Range("F19").Value = 409
$ws1.Range("F22").Value = 553
$ws1.Range("F23").Value = 3684
$ws1.Range("F25").Value = 44
$ws1.Range("F26").Value = 33
$ws1.Range("F28").Value = 3188
$ws1.Range("F30").Value = 291
$ws1.Range("G30").Value = "已售罄"
$ws1.Range("F32").Value = 363
$ws1.Range("F33").Value = 140
$ws1.Range("F34").Value = 341
$ws1.Range("F35").Value = 1015
$ws1.Range("F36").Value = 677
$ws1.Range("F38").Value = 888
$ws1.Range("F39").Value = 2670
$ws1.Range("F40").Value = 52
$ws1.Range("F43").Value = 3181
$ws1.Range("F45").Value = 2303

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 110
$ws2.Range("F3").Value = 136
$ws2.Range("F4").Value = 11
$ws2.Range("F5").Value = 56
$ws2.Range("F6").Value = 7
$ws2.Range("F10").Value = 3

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 1342

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1342
$ws4.Range("F4").Value = 19
$ws4.Range("F5").Value = 7731
$ws4.Range("F7").Value = 114
$ws4.Range("F8").Value = 2105
$ws4.Range("F9").Value = 8513
$ws4.Range("F10").Value = 11
$ws4.Range("F11").Value = 84
$ws4.Range("F12").Value = 5711
$ws4.Range("F14").Value = 2661
$ws4.Range("F15").Value = 1169
$ws4.Range("F16").Value = 4600
$ws4.Range("F17").Value = 409
$ws4.Range("F19").Value = 110
$ws4.Range("F21").Value = 136
$ws4.Range("F22").Value = 553
$ws4.Range("F23").Value = 11
$ws4.Range("F24").Value = 3684
$ws4.Range("F26").Value = 44
$ws4.Range("F27").Value = 33
$ws4.Range("F29").Value = 3189
$ws4.Range("F31").Value = 363
$ws4.Range("F32").Value = 140
$ws4.Range("F33").Value = 341
$ws4.Range("F34").Value = 56
$ws4.Range("F35").Value = 1015
$ws4.Range("F36").Value = 677
$ws4.Range("F38").Value = 888
$ws4.Range("F40").Value = 2670
$ws4.Range("F41").Value = 52
$ws4.Range("F44").Value = 3181
$ws4.Range("F46").Value = 2303
